$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("175").Insert()

$ws.Range("A175").Value = 5
$ws.Range("B175").Value = "Macroferia Regional de Talca"
$ws.Range("C175").Value = "Maule"
$ws.Range("D175").Value = 44606
$ws.Range("E175").Value = 7
$ws.Range("F175").Value = 100114013
$ws.Range("G175").Value = "Zanahoria"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 400
$ws.Range("K175").Value = 7500
$ws.Range("L175").Value = 7500
$ws.Range("M175").Value = 7500
$ws.Range("N175").Value = "$/saco 20 kilos"
$ws.Range("O175").Value = "Región de Ñuble"
$ws.Range("P175").Value = 375
$ws.Range("Q175").Value = 20
$ws.Range("R175").Value = "Hortaliza"
